$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 402 (Especial/Primera/Segunda
# group for date 44386). This shifts the existing rows 402:525 down to
# 405:528, matching the diff's row-shift pattern (dimension grows from
# A1:T525 to A1:T528).
$ws.Rows("402:404").Insert()

# Populate the 3 newly-inserted rows with the new weekly price group
# (date 44551 == 2021-12-21).

# Row 402 - Especial
$ws.Cells.Item(402, 1).Value = 8
$ws.Cells.Item(402, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(402, 3).Value = "Coquimbo"
$ws.Cells.Item(402, 4).Value = 44551
$ws.Cells.Item(402, 5).Value = 4
$ws.Cells.Item(402, 6).Value = "Fruta"
$ws.Cells.Item(402, 7).Value = 100101
$ws.Cells.Item(402, 8).Value = "Berries"
$ws.Cells.Item(402, 9).Value = 100112025
$ws.Cells.Item(402, 10).Value = "Frutilla"
$ws.Cells.Item(402, 11).Value = "Sin especificar"
$ws.Cells.Item(402, 12).Value = "Especial"
$ws.Cells.Item(402, 13).Value = 400
$ws.Cells.Item(402, 14).Value = 12000
$ws.Cells.Item(402, 15).Value = 13000
$ws.Cells.Item(402, 16).Value = 12500
$ws.Cells.Item(402, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(402, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(402, 19).Value = 1786
$ws.Cells.Item(402, 20).Value = 7

# Row 403 - Primera
$ws.Cells.Item(403, 1).Value = 8
$ws.Cells.Item(403, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(403, 3).Value = "Coquimbo"
$ws.Cells.Item(403, 4).Value = 44551
$ws.Cells.Item(403, 5).Value = 4
$ws.Cells.Item(403, 6).Value = "Fruta"
$ws.Cells.Item(403, 7).Value = 100101
$ws.Cells.Item(403, 8).Value = "Berries"
$ws.Cells.Item(403, 9).Value = 100112025
$ws.Cells.Item(403, 10).Value = "Frutilla"
$ws.Cells.Item(403, 11).Value = "Sin especificar"
$ws.Cells.Item(403, 12).Value = "Primera"
$ws.Cells.Item(403, 13).Value = 340
$ws.Cells.Item(403, 14).Value = 10000
$ws.Cells.Item(403, 15).Value = 11000
$ws.Cells.Item(403, 16).Value = 10500
$ws.Cells.Item(403, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(403, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(403, 19).Value = 1500
$ws.Cells.Item(403, 20).Value = 7

# Row 404 - Segunda
$ws.Cells.Item(404, 1).Value = 8
$ws.Cells.Item(404, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(404, 3).Value = "Coquimbo"
$ws.Cells.Item(404, 4).Value = 44551
$ws.Cells.Item(404, 5).Value = 4
$ws.Cells.Item(404, 6).Value = "Fruta"
$ws.Cells.Item(404, 7).Value = 100101
$ws.Cells.Item(404, 8).Value = "Berries"
$ws.Cells.Item(404, 9).Value = 100112025
$ws.Cells.Item(404, 10).Value = "Frutilla"
$ws.Cells.Item(404, 11).Value = "Sin especificar"
$ws.Cells.Item(404, 12).Value = "Segunda"
$ws.Cells.Item(404, 13).Value = 300
$ws.Cells.Item(404, 14).Value = 8000
$ws.Cells.Item(404, 15).Value = 9000
$ws.Cells.Item(404, 16).Value = 8500
$ws.Cells.Item(404, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(404, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(404, 19).Value = 1214
$ws.Cells.Item(404, 20).Value = 7
